$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches source inlineStr formatting)
$dCells = @("D2","D3","D5","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '30.704.04'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '1.948.58'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '247.16'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4841'
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").Value = '0.06822'
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("D10").Value = '112.55'
$ws.Range("E10").Value = '  +2.48%  '
$ws.Range("D11").Value = '19.41'
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").Value = '1.943.92'
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '0.07672'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '5.514'
$ws.Range("E14").Value = '  +4.68%  '
$ws.Range("D15").Value = '0.6903'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").Value = '296.17'
$ws.Range("E16").Value = '  +7.26%  '
$ws.Range("D17").Value = '30.765.15'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '13.29'
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").Value = '5.673'
$ws.Range("E19").Value = '  +3.53%  '
$ws.Range("D20").Value = '0.000007714'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '2.206.41'
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '6.607'
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").Value = '9.871'
$ws.Range("E25").Value = '  +4.64%  '
$ws.Range("D26").Value = '168.33'
$ws.Range("E26").Value = '  +2.84%  '
$ws.Range("D27").Value = '20.34'
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").Value = '2.193'
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("D29").Value = '0.1090'
$ws.Range("E29").Value = '  +3.63%  '
$ws.Range("D30").Value = '1.440'
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").Value = '4.741'
$ws.Range("E31").Value = '  +17.08%  '
$ws.Range("D32").Value = '4.488'
$ws.Range("E32").Value = '  +8.18%  '
$ws.Range("D33").Value = '0.05095'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").Value = '0.7780'
$ws.Range("E34").Value = '  +6.69%  '
$ws.Range("D35").Value = '1.162'
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("D36").Value = '0.02090'
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '2.700'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").Value = '2.054'
$ws.Range("D40").Value = '111.47'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").Value = '0.4469'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").Value = '0.8741'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = '5.920'
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("D44").Value = '70.17'
$ws.Range("E44").Value = '  +3.50%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '7.378'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = '9.445'
$ws.Range("E47").Value = '  +2.44%  '
$ws.Range("D48").Value = '48.65'
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("E49").Value = '  +1.37%  '
$ws.Range("D50").Value = '35.75'
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("D51").Value = '0.2533'
$ws.Range("E51").Value = '  +1.71%  '
